$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = $null

$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = $null

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null

$ws.Range("H130").Value = 96653.664
$ws.Range("J130").Value = 96653.664
$ws.Range("L130").Value = 96653.664
$ws.Range("N130").Value = -106693.664

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1814.5
$ws.Range("I97").Value = 1263.909
$ws.Range("J97").Value = 3833.3333
$ws.Range("K97").Value = 1263.909
$ws.Range("L97").Value = 3833.3333
$ws.Range("M97").Value = -767.9090000000001
$ws.Range("N97").Value = -4825.3333

$ws.Range("H106").Value = 19749.25
$ws.Range("J106").Value = 19749.25
$ws.Range("L106").Value = 19749.25
$ws.Range("N106").Value = -22273.25

$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = $null

$ws.Range("H111").Value = 62322
$ws.Range("J111").Value = 62322
$ws.Range("L111").Value = 62322
$ws.Range("N111").Value = -70502

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = $null

$ws.Range("H125").Value = 39999.5
$ws.Range("J125").Value = 39999.5
$ws.Range("L125").Value = 39999.5
$ws.Range("N125").Value = -49839.5

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = $null

$ws.Range("H94").Value = 5332.1113
$ws.Range("I94").Value = 4497.25
$ws.Range("J94").Value = 6000
$ws.Range("K94").Value = 4497.25
$ws.Range("L94").Value = 6000
$ws.Range("M94").Value = -4046.25
$ws.Range("N94").Value = -6902

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = $null

$ws.Range("H105").Value = 3527.375
$ws.Range("I105").Value = 2804
$ws.Range("K105").Value = 2804
$ws.Range("M105").Value = -1057

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = $null

$ws.Range("H115").Value = 90000
$ws.Range("I115").Value = 90000
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 90000
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -88433
$ws.Range("N115").Value = $null

$ws.Range("H119").Value = 71757.5
$ws.Range("J119").Value = 71757.5
$ws.Range("L119").Value = 71757.5
$ws.Range("N119").Value = -81433.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2438
$ws.Range("I3").Value = 984
$ws.Range("K3").Value = 984
$ws.Range("M3").Value = -871

$ws.Range("H22").Value = 1250
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -150

$ws.Range("H53").Value = 99994
$ws.Range("J53").Value = 99994
$ws.Range("L53").Value = 99994
$ws.Range("N53").Value = -101208

$ws.Range("H105").Value = 1874.75
$ws.Range("I105").Value = 999.6667
$ws.Range("J105").Value = 4500
$ws.Range("K105").Value = 999.6667
$ws.Range("L105").Value = 4500
$ws.Range("M105").Value = 747.3333
$ws.Range("N105").Value = -7994

$ws.Range("H132").Value = 3401.8
$ws.Range("I132").Value = 3401.8
$ws.Range("K132").Value = 10205.4
$ws.Range("M132").Value = -7675.400000000001

$ws.Range("H139").Value = 99994
$ws.Range("J139").Value = 99994
$ws.Range("L139").Value = 99994
$ws.Range("N139").Value = -110274

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = $null
$ws.Range("N36").Value = $null

$ws.Range("H51").Value = 968.75
$ws.Range("I51").Value = 968.75
$ws.Range("K51").Value = 2906.25
$ws.Range("M51").Value = -2446.25

$ws.Range("H118").Value = 500
$ws.Range("I118").Value = 500
$ws.Range("K118").Value = 1500
$ws.Range("M118").Value = -257

$ws.Range("H131").Value = 2161.2144
$ws.Range("J131").Value = 2408.818
$ws.Range("L131").Value = 7226.454000000001
$ws.Range("N131").Value = -17306.454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 646735.5600000001
$ws.Range("I3").Value = 627999.75
$ws.Range("J3").Value = 671716.7
$ws.Range("K3").Value = 627999.75
$ws.Range("L3").Value = 671716.7
$ws.Range("M3").Value = -627883.75
$ws.Range("N3").Value = -671948.7

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = $null

$ws.Range("H105").Value = 10000
$ws.Range("J105").Value = 10000
$ws.Range("L105").Value = 10000
$ws.Range("N105").Value = -16988

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null

$ws.Range("H132").Value = 7857.2856
$ws.Range("I132").Value = 7500.1665
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 22500.4995
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -19970.4995
$ws.Range("N132").Value = -35060

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1649.55
$ws.Range("I22").Value = 1593.1875
$ws.Range("J22").Value = 1875
$ws.Range("K22").Value = 1593.1875
$ws.Range("L22").Value = 1875
$ws.Range("M22").Value = -1298.1875
$ws.Range("N22").Value = -2465

$ws.Range("H27").Value = 1649.55
$ws.Range("I27").Value = 1593.1875
$ws.Range("J27").Value = 1875
$ws.Range("K27").Value = 1593.1875
$ws.Range("L27").Value = 1875
$ws.Range("M27").Value = -1486.1875
$ws.Range("N27").Value = -2089

$ws.Range("H55").Value = 2975.1
$ws.Range("I55").Value = 2219
$ws.Range("K55").Value = 2219
$ws.Range("M55").Value = -2046

$ws.Range("H68").Value = 3192.077
$ws.Range("I68").Value = 2749.7
$ws.Range("J68").Value = 4666.6665
$ws.Range("K68").Value = 2749.7
$ws.Range("L68").Value = 4666.6665
$ws.Range("M68").Value = -2000.7
$ws.Range("N68").Value = -6164.6665

$ws.Range("H71").Value = 3192.077
$ws.Range("I71").Value = 2749.7
$ws.Range("J71").Value = 4666.6665
$ws.Range("K71").Value = 13748.5
$ws.Range("L71").Value = 23333.3325
$ws.Range("M71").Value = -10004.5
$ws.Range("N71").Value = -30821.3325

$ws.Range("H82").Value = 1690.8948
$ws.Range("I82").Value = 1623.8
$ws.Range("J82").Value = 1765.4445
$ws.Range("K82").Value = 1623.8
$ws.Range("L82").Value = 1765.4445
$ws.Range("M82").Value = -1262.8
$ws.Range("N82").Value = -2487.4445

$ws.Range("H85").Value = 1690.8948
$ws.Range("I85").Value = 1623.8
$ws.Range("J85").Value = 1765.4445
$ws.Range("K85").Value = 1623.8
$ws.Range("L85").Value = 1765.4445
$ws.Range("M85").Value = -375.8
$ws.Range("N85").Value = -4261.4445

$ws.Range("H100").Value = 8300
$ws.Range("J100").Value = 2450
$ws.Range("L100").Value = 2450
$ws.Range("N100").Value = -3532

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = $null

$ws.Range("H124").Value = 99983
$ws.Range("J124").Value = 99983
$ws.Range("L124").Value = 99983
$ws.Range("N124").Value = -109803

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 506450.4
$ws.Range("I3").Value = 1003542
$ws.Range("J3").Value = 9358.799999999999
$ws.Range("K3").Value = 1003542
$ws.Range("L3").Value = 9358.799999999999
$ws.Range("M3").Value = -1003428
$ws.Range("N3").Value = -9586.799999999999

$ws.Range("H45").Value = 31083.334
$ws.Range("J45").Value = 30000
$ws.Range("L45").Value = 30000
$ws.Range("N45").Value = -30982

$ws.Range("H119").Value = 35000
$ws.Range("J119").Value = 35000
$ws.Range("L119").Value = 35000
$ws.Range("N119").Value = -44676

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null
